# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (positioned right before "总计"),
#    populated like the "2021-Q4" sheet but with the 2022-Q1 figures.
# 2. Insert a new top data row into "总计" summarizing the new quarter,
#    shifting the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted before "总计"
# ---------------------------------------------------------------------
$wsTotalBefore = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Add($wsTotalBefore)
$wsQ1.Name = "2022-Q1"

# Sheet references captured before the Add() call above can become
# stale once the sheet collection is restructured, so re-resolve the
# sheets we still need to touch by name.
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ1 = $wb.Worksheets.Item("2022-Q1")

# Clone the layout/values/styles of "2021-Q4" (same table shape), then
# patch in the cells that actually differ for 2022-Q1. Column A (row
# index) is copied separately so an empty A1 cell isn't materialised.
$wsQ4.Range("B1:H3").Copy($wsQ1.Range("B1"))
$wsQ4.Range("A2:A3").Copy($wsQ1.Range("A2"))

# Cells whose text values change vs. 2021-Q4. Force a text format first
# so the numeric-looking strings ("6.47", "0.0479", ...) are kept as
# text (matching 005459/005460 staying text with their leading zeros).
$changedCells = "D2", "E2", "G2", "D3", "E3", "G3"
foreach ($addr in $changedCells) {
    $wsQ1.Range($addr).NumberFormat = "@"
}

$wsQ1.Range("D2").Value = "6.47"
$wsQ1.Range("E2").Value = "39.69"
$wsQ1.Range("G2").Value = "0.0479"

$wsQ1.Range("D3").Value = "2.79"
$wsQ1.Range("E3").Value = "39.69"
$wsQ1.Range("G3").Value = "0.0206"

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert the 2022-Q1 summary row at the top of the data
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# Re-apply the row-index column style (bold/centered/bordered, same as
# the other "A" column cells) that the freshly inserted row doesn't
# carry, then clear the stray formatting Insert() propagated onto B:D.
$wsTotal.Range("A3").Copy($wsTotal.Range("A2"))
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.07000000000000001

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
